$wb = $excel.ActiveWorkbook

$wsAtrasos = $wb.Worksheets.Item("Reporte atrasos")
$wsPlanificacion = $wb.Worksheets.Item("Reporte planificación")

# --- Sheet "Reporte atrasos": update dates ---
$wsAtrasos.Range("C4").Value = 42906
$wsAtrasos.Range("C5").Value = 42920

# --- Sheet "Reporte planificación": update dates ---
# Row 4
$wsPlanificacion.Range("B4").Value = 42880
$wsPlanificacion.Range("C4").Value = 42880
$wsPlanificacion.Range("E4").Value = 42881
$wsPlanificacion.Range("F4").Value = 42881
$wsPlanificacion.Range("H4").Value = 42884
$wsPlanificacion.Range("I4").Value = 42887
$wsPlanificacion.Range("K4").Value = 42905
$wsPlanificacion.Range("L4").Value = 42906

# Row 5
$wsPlanificacion.Range("B5").Value = 42880
$wsPlanificacion.Range("C5").Value = 42880
$wsPlanificacion.Range("E5").Value = 42881
$wsPlanificacion.Range("F5").Value = 42881
$wsPlanificacion.Range("H5").Value = 42884
$wsPlanificacion.Range("I5").Value = 42888
$wsPlanificacion.Range("K5").Value = 42907
$wsPlanificacion.Range("L5").Value = 42913

# Row 6
$wsPlanificacion.Range("B6").Value = 42880
$wsPlanificacion.Range("C6").Value = 42880
$wsPlanificacion.Range("E6").Value = 42881
$wsPlanificacion.Range("F6").Value = 42881
$wsPlanificacion.Range("H6").Value = 42884
$wsPlanificacion.Range("I6").Value = 42899
$wsPlanificacion.Range("K6").Value = 42907
$wsPlanificacion.Range("L6").Value = 42920

# Row 7
$wsPlanificacion.Range("B7").Value = 42880
$wsPlanificacion.Range("C7").Value = 42880
$wsPlanificacion.Range("E7").Value = 42881
$wsPlanificacion.Range("F7").Value = 42881
$wsPlanificacion.Range("H7").Value = 42888
$wsPlanificacion.Range("I7").Value = 42895
$wsPlanificacion.Range("K7").Value = 42914
$wsPlanificacion.Range("L7").Value = 42919

# --- Update sheet view selections ---
# "Reporte atrasos": set selection to A4 (not the active sheet/tab)
$wsAtrasos.Activate()
$wsAtrasos.Range("A4").Select()

# "Reporte planificación": remove topLeftCell (scroll back to A1) and set selection to D5;
# this sheet remains the active tab.
$wsPlanificacion.Activate()
$winPlan = $excel.ActiveWindow
$winPlan.ScrollColumn = 1
$winPlan.ScrollRow = 1
$wsPlanificacion.Range("D5").Select()
